$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32: Automata for the People
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").Value = ""

# Row 74: Adhesive of Antipathy
$ws.Range("H74").Value = 2625
$ws.Range("I74").Value = 1833.3334
$ws.Range("K74").Value = 1833.3334
$ws.Range("M74").Value = -897.3334

# Row 77: It's Gonna Grow Back (L)
$ws.Range("H77").Value = 2625
$ws.Range("I77").Value = 1833.3334
$ws.Range("K77").Value = 9166.666999999999
$ws.Range("M77").Value = -4486.666999999999

# Row 80: Cleansing the Wicked Humours
$ws.Range("H80").Value = 1047.8182
$ws.Range("I80").Value = 756.5
$ws.Range("K80").Value = 2269.5
$ws.Range("M80").Value = -1271.5

# Row 83: Washing Away the Sins (L)
$ws.Range("H83").Value = 1047.8182
$ws.Range("I83").Value = 756.5
$ws.Range("K83").Value = 6808.5
$ws.Range("M83").Value = -1816.5

# Row 88: The Grave of Hemlock Groves
$ws.Range("H88").Value = 3900.8
$ws.Range("J88").Value = 4376
$ws.Range("L88").Value = 4376
$ws.Range("N88").Value = -5188

# Row 91: Dappling the Highlands (L)
$ws.Range("H91").Value = 3900.8
$ws.Range("J91").Value = 4376
$ws.Range("L91").Value = 4376
$ws.Range("N91").Value = -7184

# Row 138: All-night Crafting
$ws.Range("H138").Value = 2297.111
$ws.Range("I138").Value = 446.75
$ws.Range("K138").Value = 1340.25
$ws.Range("M138").Value = 3799.75

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 1670
$ws.Range("I2").Value = 1255
$ws.Range("K2").Value = 1255
$ws.Range("M2").Value = -1142

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 4507
$ws.Range("I61").Value = 4000
$ws.Range("K61").Value = 4000
$ws.Range("M61").Value = -3788

# Row 76: Sometimes the South Wins
$ws.Range("H76").Value = 52333.332
$ws.Range("J76").Value = 65000
$ws.Range("L76").Value = 65000
$ws.Range("N76").Value = -65676

# Row 79: The Thriller of Autumn (L)
$ws.Range("H79").Value = 52333.332
$ws.Range("J79").Value = 65000
$ws.Range("L79").Value = 65000
$ws.Range("N79").Value = -67340

# Row 88: The Mast Chance
$ws.Range("H88").Value = 2357
$ws.Range("I88").Value = 1500
$ws.Range("J88").Value = 2499.8333
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 2499.8333
$ws.Range("M88").Value = -1094
$ws.Range("N88").Value = -3311.8333

# Row 91: The Rose and the Riveter (L)
$ws.Range("H91").Value = 2357
$ws.Range("I91").Value = 1500
$ws.Range("J91").Value = 2499.8333
$ws.Range("K91").Value = 1500
$ws.Range("L91").Value = 2499.8333
$ws.Range("M91").Value = -96
$ws.Range("N91").Value = -5307.8333

# Row 110: Scheduled Maintenance
$ws.Range("H110").Value = 999.3333
$ws.Range("I110").Value = 999.5
$ws.Range("J110").Value = 999
$ws.Range("K110").Value = 999.5
$ws.Range("L110").Value = 999
$ws.Range("M110").Value = 1045.5
$ws.Range("N110").Value = -5089

# Row 116: No Scope
$ws.Range("H116").Value = 1670
$ws.Range("I116").Value = 1255
$ws.Range("K116").Value = 1255
$ws.Range("M116").Value = 1039

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 4507
$ws.Range("I136").Value = 4000
$ws.Range("K136").Value = 12000
$ws.Range("M136").Value = -9450

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value = 1670
$ws.Range("I3").Value = 1255
$ws.Range("K3").Value = 1255
$ws.Range("M3").Value = -1141

$ws = $wb.Worksheets.Item("CRP")
# Row 5: Bowing Out
$ws.Range("H5").Value = 193.28572
$ws.Range("I5").Value = 170.6
$ws.Range("K5").Value = 170.6
$ws.Range("M5").Value = -58.59999999999999

# Row 25: Bowing to Necessity
$ws.Range("H25").Value = 11
$ws.Range("I25").Value = 11
$ws.Range("K25").Value = 11
$ws.Range("M25").Value = 163

# Row 41: The Lone Bowman
$ws.Range("H41").Value = 5000
$ws.Range("I41").Value = 5000
$ws.Range("K41").Value = 5000
$ws.Range("M41").Value = -4572

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 9874.25
$ws.Range("I58").Value = 9874.25
$ws.Range("K58").Value = 9874.25
$ws.Range("M58").Value = -9671.25

# Row 62: Splinter in the Sewers
$ws.Range("H62").Value = 6138
$ws.Range("I62").Value = 4209.5
$ws.Range("K62").Value = 4209.5
$ws.Range("M62").Value = -3585.5

# Row 65: The Lumber of Their Discontent (L)
$ws.Range("H65").Value = 6138
$ws.Range("I65").Value = 4209.5
$ws.Range("K65").Value = 21047.5
$ws.Range("M65").Value = -17927.5

# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 6327.5713
$ws.Range("I122").Value = 6882.1665
$ws.Range("K122").Value = 20646.4995
$ws.Range("M122").Value = -18196.4995

# Row 136: Turali Quality
$ws.Range("H136").Value = 9874.25
$ws.Range("I136").Value = 9874.25
$ws.Range("K136").Value = 29622.75
$ws.Range("M136").Value = -27072.75

$ws = $wb.Worksheets.Item("CUL")
# Row 64: The Aroma of Faith
$ws.Range("H64").Value = 4500
$ws.Range("J64").Value = 4500
$ws.Range("L64").Value = 13500
$ws.Range("N64").Value = -14040

# Row 67: Soup's On (L)
$ws.Range("H67").Value = 4500
$ws.Range("J67").Value = 4500
$ws.Range("L67").Value = 13500
$ws.Range("N67").Value = -15372

$ws = $wb.Worksheets.Item("GSM")
# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 7494.8
$ws.Range("I122").Value = 2269
$ws.Range("J122").Value = 15333.5
$ws.Range("K122").Value = 6807
$ws.Range("L122").Value = 46000.5
$ws.Range("M122").Value = -4357
$ws.Range("N122").Value = -50900.5

$ws = $wb.Worksheets.Item("LTW")
# Row 46: Supply Side Logic
$ws.Range("H46").Value = 7900
$ws.Range("I46").Value = 7900
$ws.Range("K46").Value = 7900
$ws.Range("M46").Value = -7712

# Row 55: It's Not a Job, It's a Calling
$ws.Range("H55").Value = 2483.7
$ws.Range("I55").Value = 2367.125
$ws.Range("J55").Value = 2950
$ws.Range("K55").Value = 2367.125
$ws.Range("L55").Value = 2950
$ws.Range("M55").Value = -2194.125
$ws.Range("N55").Value = -3296

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 8716.6875
$ws.Range("I132").Value = 8580.333000000001
$ws.Range("J132").Value = 9125.75
$ws.Range("K132").Value = 25740.999
$ws.Range("L132").Value = 27377.25
$ws.Range("M132").Value = -23210.999
$ws.Range("N132").Value = -32437.25

$ws = $wb.Worksheets.Item("WVR")
# Row 122: Heavy Armoire
$ws.Range("H122").Value = 3998.6
$ws.Range("I122").Value = 3497
$ws.Range("J122").Value = 4333
$ws.Range("K122").Value = 10491
$ws.Range("L122").Value = 12999
$ws.Range("M122").Value = -8041
$ws.Range("N122").Value = -17899

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 3148.5
$ws.Range("I126").Value = 3248.125
$ws.Range("K126").Value = 9744.375
$ws.Range("M126").Value = -7274.375

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 7179.5
$ws.Range("I132").Value = 5739.5
$ws.Range("K132").Value = 17218.5
$ws.Range("M132").Value = -14688.5
